$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows that currently hold a student record with a GitHub hyperlink in
# column C (every data row except the two blank separator rows 4 and 17).
$rows = @(2,3,5,6,7,8,9,10,11,12,13,14,15,16,18,19,20,21,22,23,24,25,26,27,28,29,30,31)

# Remember the hyperlink target for every one of those rows (the display
# text in column C is the URL itself) before we start mutating the sheet.
$urls = @{}
foreach ($r in $rows) {
    $urls[$r] = $ws.Range("C$r").Value2
}

# Capture the direct cell formatting used throughout column C so it can be
# restored after re-creating the hyperlinks (Hyperlinks.Add forces the
# built-in "Hyperlink" look on the cell, which is not what the original
# file uses - column C cells already carry their own direct font/fill).
$sample = $ws.Range("C4")
$cFontName = $sample.Font.Name
$cFontSize = $sample.Font.Size
$cFontUnderline = $sample.Font.Underline
$cFontColor = $sample.Font.Color
$cInteriorColor = $sample.Interior.Color
$cNumberFormat = $sample.NumberFormat

# Remove the two students who left the roster: Camelia Ignat (row 6) and
# Magdalena Mostazo (row 18). Use ClearContents (not a full row delete) so
# the rows collapse into blank placeholder rows exactly like rows 4/17,
# keeping every other row number, style and hyperlink anchor unchanged.
$ws.Range("A6:D6").ClearContents()
$ws.Range("A18:D18").ClearContents()

# The hyperlink collection in this runtime can only be cleared in bulk, so
# drop every hyperlink and re-add the ones that should remain. This also
# renumbers the relationship ids sequentially, just like a normal Excel
# re-save after removing two hyperlinks from the middle of the list.
$ws.Hyperlinks.Delete()

foreach ($r in $rows) {
    if ($r -eq 6 -or $r -eq 18) { continue }
    $url = $urls[$r]
    $cell = $ws.Range("C$r")
    $ws.Hyperlinks.Add($cell, $url, [Type]::Missing, [Type]::Missing, $url) | Out-Null

    # Restore the original direct formatting that Hyperlinks.Add overwrote.
    $cell.Font.Name = $cFontName
    $cell.Font.Size = $cFontSize
    $cell.Font.Underline = $cFontUnderline
    $cell.Font.Color = $cFontColor
    $cell.Interior.Color = $cInteriorColor
    $cell.NumberFormat = $cNumberFormat
}

# Reflect the edit in the active selection: the last thing touched was the
# newly blanked A18:D18 row.
$ws.Range("A18:D18").Select()
